$d = $word.ActiveDocument

# 1) TRADE table: "user serial references USERS (user_id) NOT NULL,"
#    -> "client serial references USERS (user_id) NOT NULL,"
$d.Content.Find.Execute(
    "user serial references USERS (user_id) NOT NULL,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "client serial references USERS (user_id) NOT NULL,", 2) | Out-Null

# 2) ISSUES table: "user_id integer references users (user_id) NOT NULL,"
#    -> "client serial references USERS (user_id) NOT NULL,"
$d.Content.Find.Execute(
    "user_id integer references users (user_id) NOT NULL,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "client serial references USERS (user_id) NOT NULL,", 2) | Out-Null

# 3) ISSUES table: "property_id integer references properties (properties_id) NOT NULL"
#    -> "property serial references PROPERTIES (properties_id) NOT NULL"
$d.Content.Find.Execute(
    "property_id integer references properties (properties_id) NOT NULL",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "property serial references PROPERTIES (properties_id) NOT NULL", 2) | Out-Null
